# Applies the "Trade #27 closed" update + new open trade (#88) to the
# live trading results workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.57   # Current Capital
$summary.Range("B4").Value = 0.37      # Total P&L $
$summary.Range("B6").Value = 55        # Total Trades
$summary.Range("B7").Value = 27        # Winning Trades
$summary.Range("B9").Value = 49.09     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.57
$status.Range("D5").Value = 22
$status.Range("E5").Value = 0.26
$status.Range("F5").Value = 0.57
$status.Range("G5").Value = 59.09

# ---------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Trade #55 (row 56) closes out with an early exit.
$allTrades.Range("G56").Value = 0.92
$allTrades.Range("H56").Value = "CLOSED"
$allTrades.Range("I56").Value = 2.2222
$allTrades.Range("J56").Value = 0.02
$allTrades.Range("K56").Value = 100.57
$allTrades.Range("L56").Value = "early_exit"
$allTrades.Range("M56").Value = 0.14

# New trade #88 (row 89) opened.
$allTrades.Range("A89").Value = 88
$allTrades.Range("B89").Value = "'2026-02-17"
$allTrades.Range("C89").Value = "'20:54:14"
$allTrades.Range("D89").Value = "MarketMaking"
$allTrades.Range("E89").Value = "DOWN"
$allTrades.Range("F89").Value = 0.9
$allTrades.Range("H89").Value = "OPEN"
$allTrades.Range("I89").Value = 0
$allTrades.Range("J89").Value = 0
$allTrades.Range("K89").Value = 100.5534535840667
$allTrades.Range("M89").Value = 0
$allTrades.Range("N89").Value = 0
$allTrades.Range("O89").Value = 0
$allTrades.Range("P89").Value = 0.6
$allTrades.Range("Q89").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------
# MarketMaking sheet
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

# Trade #55 (row 23) closes out with an early exit.
$mm.Range("G23").Value = 0.92
$mm.Range("H23").Value = "CLOSED"
$mm.Range("I23").Value = 2.2222
$mm.Range("J23").Value = 0.02
$mm.Range("K23").Value = 100.57
$mm.Range("P23").Value = "early_exit"
$mm.Range("Q23").Value = 0.14

# New trade #88 (row 56) opened.
$mm.Range("A56").Value = 88
$mm.Range("B56").Value = "'2026-02-17"
$mm.Range("C56").Value = "'20:54:14"
$mm.Range("D56").Value = "MarketMaking"
$mm.Range("E56").Value = "DOWN"
$mm.Range("F56").Value = 0.9
$mm.Range("H56").Value = "OPEN"
$mm.Range("I56").Value = 0
$mm.Range("J56").Value = 0
$mm.Range("K56").Value = 100.5534535840667
$mm.Range("L56").Value = 0
$mm.Range("M56").Value = 0
$mm.Range("N56").Value = 0.6
$mm.Range("O56").Value = "Normal spread capture: 19600 bps"
$mm.Range("Q56").Value = 0
